$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds "Municipio" names; homologate to INEGI "CVE_MUN" codes.
# Header (row 1): "Municipio" -> "CVE_MUN"
$ws.Cells.Item(1, 1).Value = "CVE_MUN"

# Rows 2-85: municipality name -> numeric CVE_MUN code.
# The code must be written as TEXT (it is a key, not a quantity), matching
# the rest of the sheet, which stores every value as a shared text string.
# Pre-setting NumberFormat to "@" (Text) makes Excel keep the digit string
# as text instead of silently parsing it into a Number; ClearFormats()
# afterwards drops that temporary number-format again so no cell ends up
# with a different style than before.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "13001"  # Acatlán
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "13002"  # Acaxochitlán
$ws.Cells.Item(3, 1).ClearFormats()
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "13003"  # Actopan
$ws.Cells.Item(4, 1).ClearFormats()
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "13004"  # Agua Blanca de Iturbide
$ws.Cells.Item(5, 1).ClearFormats()
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "13005"  # Ajacuba
$ws.Cells.Item(6, 1).ClearFormats()
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "13006"  # Alfajayucan
$ws.Cells.Item(7, 1).ClearFormats()
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "13007"  # Almoloya
$ws.Cells.Item(8, 1).ClearFormats()
$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "13008"  # Apan
$ws.Cells.Item(9, 1).ClearFormats()
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "13010"  # Atitalaquia
$ws.Cells.Item(10, 1).ClearFormats()
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "13011"  # Atlapexco
$ws.Cells.Item(11, 1).ClearFormats()
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "13013"  # Atotonilco de Tula
$ws.Cells.Item(12, 1).ClearFormats()
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "13012"  # Atotonilco el Grande
$ws.Cells.Item(13, 1).ClearFormats()
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = "13014"  # Calnali
$ws.Cells.Item(14, 1).ClearFormats()
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "13015"  # Cardonal
$ws.Cells.Item(15, 1).ClearFormats()
$ws.Cells.Item(16, 1).NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = "13017"  # Chapantongo
$ws.Cells.Item(16, 1).ClearFormats()
$ws.Cells.Item(17, 1).NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = "13018"  # Chapulhuacán
$ws.Cells.Item(17, 1).ClearFormats()
$ws.Cells.Item(18, 1).NumberFormat = "@"
$ws.Cells.Item(18, 1).Value = "13019"  # Chilcuautla
$ws.Cells.Item(18, 1).ClearFormats()
$ws.Cells.Item(19, 1).NumberFormat = "@"
$ws.Cells.Item(19, 1).Value = "13016"  # Cuautepec de Hinojosa
$ws.Cells.Item(19, 1).ClearFormats()
$ws.Cells.Item(20, 1).NumberFormat = "@"
$ws.Cells.Item(20, 1).Value = "13009"  # El Arenal
$ws.Cells.Item(20, 1).ClearFormats()
$ws.Cells.Item(21, 1).NumberFormat = "@"
$ws.Cells.Item(21, 1).Value = "13020"  # Eloxochitlán
$ws.Cells.Item(21, 1).ClearFormats()
$ws.Cells.Item(22, 1).NumberFormat = "@"
$ws.Cells.Item(22, 1).Value = "13021"  # Emiliano Zapata
$ws.Cells.Item(22, 1).ClearFormats()
$ws.Cells.Item(23, 1).NumberFormat = "@"
$ws.Cells.Item(23, 1).Value = "13022"  # Epazoyucan
$ws.Cells.Item(23, 1).ClearFormats()
$ws.Cells.Item(24, 1).NumberFormat = "@"
$ws.Cells.Item(24, 1).Value = "13023"  # Francisco I. Madero
$ws.Cells.Item(24, 1).ClearFormats()
$ws.Cells.Item(25, 1).NumberFormat = "@"
$ws.Cells.Item(25, 1).Value = "13024"  # Huasca de Ocampo
$ws.Cells.Item(25, 1).ClearFormats()
$ws.Cells.Item(26, 1).NumberFormat = "@"
$ws.Cells.Item(26, 1).Value = "13025"  # Huautla
$ws.Cells.Item(26, 1).ClearFormats()
$ws.Cells.Item(27, 1).NumberFormat = "@"
$ws.Cells.Item(27, 1).Value = "13026"  # Huazalingo
$ws.Cells.Item(27, 1).ClearFormats()
$ws.Cells.Item(28, 1).NumberFormat = "@"
$ws.Cells.Item(28, 1).Value = "13027"  # Huehuetla
$ws.Cells.Item(28, 1).ClearFormats()
$ws.Cells.Item(29, 1).NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = "13028"  # Huejutla de Reyes
$ws.Cells.Item(29, 1).ClearFormats()
$ws.Cells.Item(30, 1).NumberFormat = "@"
$ws.Cells.Item(30, 1).Value = "13029"  # Huichapan
$ws.Cells.Item(30, 1).ClearFormats()
$ws.Cells.Item(31, 1).NumberFormat = "@"
$ws.Cells.Item(31, 1).Value = "13030"  # Ixmiquilpan
$ws.Cells.Item(31, 1).ClearFormats()
$ws.Cells.Item(32, 1).NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "13031"  # Jacala de Ledezma
$ws.Cells.Item(32, 1).ClearFormats()
$ws.Cells.Item(33, 1).NumberFormat = "@"
$ws.Cells.Item(33, 1).Value = "13032"  # Jaltocán
$ws.Cells.Item(33, 1).ClearFormats()
$ws.Cells.Item(34, 1).NumberFormat = "@"
$ws.Cells.Item(34, 1).Value = "13033"  # Juárez Hidalgo
$ws.Cells.Item(34, 1).ClearFormats()
$ws.Cells.Item(35, 1).NumberFormat = "@"
$ws.Cells.Item(35, 1).Value = "13040"  # La Misión
$ws.Cells.Item(35, 1).ClearFormats()
$ws.Cells.Item(36, 1).NumberFormat = "@"
$ws.Cells.Item(36, 1).Value = "13034"  # Lolotla
$ws.Cells.Item(36, 1).ClearFormats()
$ws.Cells.Item(37, 1).NumberFormat = "@"
$ws.Cells.Item(37, 1).Value = "13035"  # Metepec
$ws.Cells.Item(37, 1).ClearFormats()
$ws.Cells.Item(38, 1).NumberFormat = "@"
$ws.Cells.Item(38, 1).Value = "13037"  # Metztitlán
$ws.Cells.Item(38, 1).ClearFormats()
$ws.Cells.Item(39, 1).NumberFormat = "@"
$ws.Cells.Item(39, 1).Value = "13051"  # Mineral de la Reforma
$ws.Cells.Item(39, 1).ClearFormats()
$ws.Cells.Item(40, 1).NumberFormat = "@"
$ws.Cells.Item(40, 1).Value = "13038"  # Mineral del Chico
$ws.Cells.Item(40, 1).ClearFormats()
$ws.Cells.Item(41, 1).NumberFormat = "@"
$ws.Cells.Item(41, 1).Value = "13039"  # Mineral del Monte
$ws.Cells.Item(41, 1).ClearFormats()
$ws.Cells.Item(42, 1).NumberFormat = "@"
$ws.Cells.Item(42, 1).Value = "13041"  # Mixquiahuala de Juárez
$ws.Cells.Item(42, 1).ClearFormats()
$ws.Cells.Item(43, 1).NumberFormat = "@"
$ws.Cells.Item(43, 1).Value = "13042"  # Molango de Escamilla
$ws.Cells.Item(43, 1).ClearFormats()
$ws.Cells.Item(44, 1).NumberFormat = "@"
$ws.Cells.Item(44, 1).Value = "13043"  # Nicolás Flores
$ws.Cells.Item(44, 1).ClearFormats()
$ws.Cells.Item(45, 1).NumberFormat = "@"
$ws.Cells.Item(45, 1).Value = "13044"  # Nopala de Villagrán
$ws.Cells.Item(45, 1).ClearFormats()
$ws.Cells.Item(46, 1).NumberFormat = "@"
$ws.Cells.Item(46, 1).Value = "13045"  # Omitlán de Juárez
$ws.Cells.Item(46, 1).ClearFormats()
$ws.Cells.Item(47, 1).NumberFormat = "@"
$ws.Cells.Item(47, 1).Value = "13048"  # Pachuca de Soto
$ws.Cells.Item(47, 1).ClearFormats()
$ws.Cells.Item(48, 1).NumberFormat = "@"
$ws.Cells.Item(48, 1).Value = "13047"  # Pacula
$ws.Cells.Item(48, 1).ClearFormats()
$ws.Cells.Item(49, 1).NumberFormat = "@"
$ws.Cells.Item(49, 1).Value = "13049"  # Pisaflores
$ws.Cells.Item(49, 1).ClearFormats()
$ws.Cells.Item(50, 1).NumberFormat = "@"
$ws.Cells.Item(50, 1).Value = "13050"  # Progreso de Obregón
$ws.Cells.Item(50, 1).ClearFormats()
$ws.Cells.Item(51, 1).NumberFormat = "@"
$ws.Cells.Item(51, 1).Value = "13036"  # San Agustín Metzquititlán
$ws.Cells.Item(51, 1).ClearFormats()
$ws.Cells.Item(52, 1).NumberFormat = "@"
$ws.Cells.Item(52, 1).Value = "13052"  # San Agustín Tlaxiaca
$ws.Cells.Item(52, 1).ClearFormats()
$ws.Cells.Item(53, 1).NumberFormat = "@"
$ws.Cells.Item(53, 1).Value = "13053"  # San Bartolo Tutotepec
$ws.Cells.Item(53, 1).ClearFormats()
$ws.Cells.Item(54, 1).NumberFormat = "@"
$ws.Cells.Item(54, 1).Value = "13046"  # San Felipe Orizatlán
$ws.Cells.Item(54, 1).ClearFormats()
$ws.Cells.Item(55, 1).NumberFormat = "@"
$ws.Cells.Item(55, 1).Value = "13054"  # San Salvador
$ws.Cells.Item(55, 1).ClearFormats()
$ws.Cells.Item(56, 1).NumberFormat = "@"
$ws.Cells.Item(56, 1).Value = "13055"  # Santiago de Anaya
$ws.Cells.Item(56, 1).ClearFormats()
$ws.Cells.Item(57, 1).NumberFormat = "@"
$ws.Cells.Item(57, 1).Value = "13056"  # Santiago Tulantepec de Lugo Guerrero
$ws.Cells.Item(57, 1).ClearFormats()
$ws.Cells.Item(58, 1).NumberFormat = "@"
$ws.Cells.Item(58, 1).Value = "13057"  # Singuilucan
$ws.Cells.Item(58, 1).ClearFormats()
$ws.Cells.Item(59, 1).NumberFormat = "@"
$ws.Cells.Item(59, 1).Value = "13058"  # Tasquillo
$ws.Cells.Item(59, 1).ClearFormats()
$ws.Cells.Item(60, 1).NumberFormat = "@"
$ws.Cells.Item(60, 1).Value = "13059"  # Tecozautla
$ws.Cells.Item(60, 1).ClearFormats()
$ws.Cells.Item(61, 1).NumberFormat = "@"
$ws.Cells.Item(61, 1).Value = "13060"  # Tenango de Doria
$ws.Cells.Item(61, 1).ClearFormats()
$ws.Cells.Item(62, 1).NumberFormat = "@"
$ws.Cells.Item(62, 1).Value = "13061"  # Tepeapulco
$ws.Cells.Item(62, 1).ClearFormats()
$ws.Cells.Item(63, 1).NumberFormat = "@"
$ws.Cells.Item(63, 1).Value = "13062"  # Tepehuacán de Guerrero
$ws.Cells.Item(63, 1).ClearFormats()
$ws.Cells.Item(64, 1).NumberFormat = "@"
$ws.Cells.Item(64, 1).Value = "13063"  # Tepeji del Río de Ocampo
$ws.Cells.Item(64, 1).ClearFormats()
$ws.Cells.Item(65, 1).NumberFormat = "@"
$ws.Cells.Item(65, 1).Value = "13064"  # Tepetitlán
$ws.Cells.Item(65, 1).ClearFormats()
$ws.Cells.Item(66, 1).NumberFormat = "@"
$ws.Cells.Item(66, 1).Value = "13065"  # Tetepango
$ws.Cells.Item(66, 1).ClearFormats()
$ws.Cells.Item(67, 1).NumberFormat = "@"
$ws.Cells.Item(67, 1).Value = "13067"  # Tezontepec de Aldama
$ws.Cells.Item(67, 1).ClearFormats()
$ws.Cells.Item(68, 1).NumberFormat = "@"
$ws.Cells.Item(68, 1).Value = "13068"  # Tianguistengo
$ws.Cells.Item(68, 1).ClearFormats()
$ws.Cells.Item(69, 1).NumberFormat = "@"
$ws.Cells.Item(69, 1).Value = "13069"  # Tizayuca
$ws.Cells.Item(69, 1).ClearFormats()
$ws.Cells.Item(70, 1).NumberFormat = "@"
$ws.Cells.Item(70, 1).Value = "13070"  # Tlahuelilpan
$ws.Cells.Item(70, 1).ClearFormats()
$ws.Cells.Item(71, 1).NumberFormat = "@"
$ws.Cells.Item(71, 1).Value = "13071"  # Tlahuiltepa
$ws.Cells.Item(71, 1).ClearFormats()
$ws.Cells.Item(72, 1).NumberFormat = "@"
$ws.Cells.Item(72, 1).Value = "13072"  # Tlanalapa
$ws.Cells.Item(72, 1).ClearFormats()
$ws.Cells.Item(73, 1).NumberFormat = "@"
$ws.Cells.Item(73, 1).Value = "13073"  # Tlanchinol
$ws.Cells.Item(73, 1).ClearFormats()
$ws.Cells.Item(74, 1).NumberFormat = "@"
$ws.Cells.Item(74, 1).Value = "13074"  # Tlaxcoapan
$ws.Cells.Item(74, 1).ClearFormats()
$ws.Cells.Item(75, 1).NumberFormat = "@"
$ws.Cells.Item(75, 1).Value = "13075"  # Tolcayuca
$ws.Cells.Item(75, 1).ClearFormats()
$ws.Cells.Item(76, 1).NumberFormat = "@"
$ws.Cells.Item(76, 1).Value = "13076"  # Tula de Allende
$ws.Cells.Item(76, 1).ClearFormats()
$ws.Cells.Item(77, 1).NumberFormat = "@"
$ws.Cells.Item(77, 1).Value = "13077"  # Tulancingo de Bravo
$ws.Cells.Item(77, 1).ClearFormats()
$ws.Cells.Item(78, 1).NumberFormat = "@"
$ws.Cells.Item(78, 1).Value = "13066"  # Villa de Tezontepec
$ws.Cells.Item(78, 1).ClearFormats()
$ws.Cells.Item(79, 1).NumberFormat = "@"
$ws.Cells.Item(79, 1).Value = "13078"  # Xochiatipan
$ws.Cells.Item(79, 1).ClearFormats()
$ws.Cells.Item(80, 1).NumberFormat = "@"
$ws.Cells.Item(80, 1).Value = "13079"  # Xochicoatlán
$ws.Cells.Item(80, 1).ClearFormats()
$ws.Cells.Item(81, 1).NumberFormat = "@"
$ws.Cells.Item(81, 1).Value = "13080"  # Yahualica
$ws.Cells.Item(81, 1).ClearFormats()
$ws.Cells.Item(82, 1).NumberFormat = "@"
$ws.Cells.Item(82, 1).Value = "13081"  # Zacualtipán de Ángeles
$ws.Cells.Item(82, 1).ClearFormats()
$ws.Cells.Item(83, 1).NumberFormat = "@"
$ws.Cells.Item(83, 1).Value = "13082"  # Zapotlán de Juárez
$ws.Cells.Item(83, 1).ClearFormats()
$ws.Cells.Item(84, 1).NumberFormat = "@"
$ws.Cells.Item(84, 1).Value = "13083"  # Zempoala
$ws.Cells.Item(84, 1).ClearFormats()
$ws.Cells.Item(85, 1).NumberFormat = "@"
$ws.Cells.Item(85, 1).Value = "13084"  # Zimapán
$ws.Cells.Item(85, 1).ClearFormats()
